$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    This shared string is referenced from the Overview sheet (E/F cols)
#    and from the per-locale Status column (C) on zh-cn / de-de.
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2) Populate the handback columns (I: Latest Target File, J: Latest Handback
#    File, K: Latest Handback DateTime) for both locale sheets, rows 2 & 3,
#    and add the matching hyperlink on column I (same target as column A's
#    hyperlink for that row).
# ---------------------------------------------------------------------------
$mdUrlRow2 = "https://github.com/OpenLocalizationTestOrg/oltest/blob/9c4259e484e281ef66c80b76ead61cf819be3942/e2e/e3079e2f-56cf-425c-8724-1783cce8b10b.md"
$mdUrlRow3 = "https://github.com/OpenLocalizationTestOrg/oltest/blob/9c4259e484e281ef66c80b76ead61cf819be3942/e2e/ffff276d57ca-c1b2-439c-8fc3-a2c4f13d4026.md"
$mdDisplayRow2 = "e3079e2f-56cf-425c-8724-1783cce8b10b.md"
$mdDisplayRow3 = "ffff276d57ca-c1b2-439c-8fc3-a2c4f13d4026.md"

function Set-HandbackColumns($ws, $targetFileName, $handbackDateRow2, $handbackDateRow3) {
    $ws.Range("I2").Value = $mdDisplayRow2
    $ws.Range("J2").Value = $targetFileName
    $ws.Range("K2").Value = $handbackDateRow2

    $ws.Range("I3").Value = $mdDisplayRow2
    $ws.Range("J3").Value = $targetFileName
    $ws.Range("K3").Value = $handbackDateRow3

    # Rebuild the hyperlinks collection in row order (A2, I2, A3, I3) so the
    # relationship ids come out interleaved the way Excel lays them out when
    # the new hyperlinks are inserted alongside the pre-existing ones.
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $mdUrlRow2, [Type]::Missing, [Type]::Missing, $mdDisplayRow2) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("I2"), $mdUrlRow2, [Type]::Missing, [Type]::Missing, $mdDisplayRow2) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("A3"), $mdUrlRow3, [Type]::Missing, [Type]::Missing, $mdDisplayRow3) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("I3"), $mdUrlRow2, [Type]::Missing, [Type]::Missing, $mdDisplayRow2) | Out-Null

    # Match the existing custom hyperlink look (cornflower blue, underlined)
    # instead of the engine's freshly-minted theme-hyperlink style.
    foreach ($addr in @("A2", "I2", "A3", "I3")) {
        $c = $ws.Range($addr)
        $c.Font.Underline = 2
        $c.Font.Color = 15570276
        $c.Font.Name = "Calibri"
    }
}

Set-HandbackColumns $wsZhCn "e3079e2f-56cf-425c-8724-1783cce8b10b.465a685268559761390e72b4e7b7da9986b7100f.zh-cn.xlf" "2016-08-13 11:17:56" "2016-08-13 11:17:56"
Set-HandbackColumns $wsDeDe "e3079e2f-56cf-425c-8724-1783cce8b10b.465a685268559761390e72b4e7b7da9986b7100f.de-de.xlf" "2016-08-13 11:18:08" "2016-08-13 11:18:08"

# ---------------------------------------------------------------------------
# 3) Widen columns to fit the new, longer text.
# ---------------------------------------------------------------------------
$wsOverview.Columns("E").ColumnWidth = 29.9777047293527
$wsOverview.Columns("F").ColumnWidth = 29.9777047293527

foreach ($ws in @($wsZhCn, $wsDeDe)) {
    $ws.Columns("C").ColumnWidth = 29.9777047293527
    $ws.Columns("I").ColumnWidth = 40
    $ws.Columns("J").ColumnWidth = 40
}
